# chore: adapt column header formatting to respective input file names (#7)
#
# The sheet compares an "old" (FV2404) AHB export against a "new" (FV2410)
# export. Until now the header row used generic "_old"/"_new" suffixes;
# this adapts them to the concrete format-version names, freezes the header
# row, and turns the used range into a proper Excel Table so the headers
# double as AutoFilter buttons.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row ------------------------------------------
# Columns A..J: "<Name>_old"  -> "<Name>_FV2404"
# Column  K   : "diff"        -> unchanged
# Columns L..U: "<Name>_new"  -> "<Name>_FV2410"
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($col = 1; $col -le $headers.Count; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- 2. Freeze the header row --------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap the used range in an Excel Table -----------------------------
$tableRange = $ws.Range("A1:U80")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"
